$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 956, shifting the existing data (rows 956-1010)
# down to rows 958-1012.
$ws.Rows.Item(956).Resize(2).Insert()

# Populate the two newly-inserted rows (956 and 957) with new price-report entries.
$data = New-Object 'object[,]' 2,18

# Row 956: Zafiro rojo, Primera
$data[0,0]  = 10
$data[0,1]  = "Vega Modelo de Temuco"
$data[0,2]  = "La Araucanía"
$data[0,3]  = 44516
$data[0,4]  = 9
$data[0,5]  = 100112002
$data[0,6]  = "Pimiento"
$data[0,7]  = "Zafiro rojo"
$data[0,8]  = "Primera"
$data[0,9]  = 120
$data[0,10] = 45000
$data[0,11] = 48000
$data[0,12] = 46375
$data[0,13] = "`$/caja 15 kilos"
$data[0,14] = "Región de Arica y Parinacota"
$data[0,15] = 3092
$data[0,16] = 15
$data[0,17] = "Hortaliza"

# Row 957: Zafiro verde, Primera
$data[1,0]  = 10
$data[1,1]  = "Vega Modelo de Temuco"
$data[1,2]  = "La Araucanía"
$data[1,3]  = 44516
$data[1,4]  = 9
$data[1,5]  = 100112002
$data[1,6]  = "Pimiento"
$data[1,7]  = "Zafiro verde"
$data[1,8]  = "Primera"
$data[1,9]  = 125
$data[1,10] = 35000
$data[1,11] = 35000
$data[1,12] = 35000
$data[1,13] = "`$/caja 15 kilos"
$data[1,14] = "Región de Arica y Parinacota"
$data[1,15] = 2333
$data[1,16] = 15
$data[1,17] = "Hortaliza"

$ws.Range("A956:R957").Value = $data

# Keep the date columns formatted the same way as the surrounding rows.
$ws.Range("D956:D957").NumberFormat = $ws.Range("D958").NumberFormat
